$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (tab title) for the new date
$ws.Name = "Through 2022-10-13"

# Update the row label for October to reflect the new "through" date
$ws.Range("A11").Value = "October (through 10-13)"

# Update October row (row 11) values
$ws.Range("C11").Value = 19
$ws.Range("D11").Value = 21
$ws.Range("E11").Value = 32
$ws.Range("F11").Value = 17
$ws.Range("G11").Value = 63
$ws.Range("H11").Value = 81
$ws.Range("I11").Value = 44

# Update Total row (row 12) values
$ws.Range("C12").Value = 448
$ws.Range("D12").Value = 648
$ws.Range("E12").Value = 580
$ws.Range("F12").Value = 439
$ws.Range("G12").Value = 964
$ws.Range("H12").Value = 1328
$ws.Range("I12").Value = 1322
